# Update the "想去人数" (want-to-go headcount) column F figures on each
# sheet to match the freshly generated gh-pages data snapshot.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F3").Value = 14765
$ws.Range("F5").Value = 1671
$ws.Range("F8").Value = 1313
$ws.Range("F9").Value = 1996
$ws.Range("F11").Value = 46
$ws.Range("F12").Value = 2365
$ws.Range("F14").Value = 838
$ws.Range("F15").Value = 3690
$ws.Range("F18").Value = 2739
$ws.Range("F19").Value = 710
$ws.Range("F22").Value = 1947
$ws.Range("F24").Value = 1669
$ws.Range("F27").Value = 7699
$ws.Range("F28").Value = 5308
$ws.Range("F29").Value = 337
$ws.Range("F31").Value = 731
$ws.Range("F32").Value = 738
$ws.Range("F33").Value = 3423
$ws.Range("F39").Value = 4528
$ws.Range("F40").Value = 757
$ws.Range("F42").Value = 358

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F18").Value = 129

$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F2").Value = 8081
$ws.Range("F4").Value = 1164

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value = 8081
$ws.Range("F5").Value = 1164
$ws.Range("F6").Value = 14765
$ws.Range("F9").Value = 1671
$ws.Range("F11").Value = 1313
$ws.Range("F12").Value = 1996
$ws.Range("F14").Value = 46
$ws.Range("F17").Value = 3690
$ws.Range("F19").Value = 2739
$ws.Range("F20").Value = 710
$ws.Range("F23").Value = 1947
$ws.Range("F29").Value = 1669
$ws.Range("F33").Value = 7699
$ws.Range("F34").Value = 5308
$ws.Range("F35").Value = 337
$ws.Range("F36").Value = 731
$ws.Range("F37").Value = 3423
$ws.Range("F43").Value = 4528
$ws.Range("F44").Value = 757
$ws.Range("F46").Value = 358

